# Auto-generated edit script for cryptos.xlsx update
# Updates Coin/Link/Price/Volume(1h) columns for rows 2-51 per the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.177.73'
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").Value = '2.484.18'
$ws.Range("E3").Value = '  +0.68%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = "'519.28"
$ws.Range("E5").Value = '  +0.27%  '
$ws.Range("D6").Value = "'131.99"
$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = '  -0.46%  '
$ws.Range("E8").Value = '  -0.61%  '
$ws.Range("D9").Value = '2.517.65'
$ws.Range("E9").Value = '  +1.97%  '
$ws.Range("E10").Value = '  -1.67%  '
$ws.Range("D11").Value = "'0.156"
$ws.Range("E11").Value = '  -0.11%  '
$ws.Range("D12").Value = "'5.19"
$ws.Range("E12").Value = '  -2.61%  '
$ws.Range("E13").Value = '  -2.41%  '
$ws.Range("D14").Value = '2.933.59'
$ws.Range("E14").Value = '  +0.90%  '
$ws.Range("D15").Value = '58.089.20'
$ws.Range("E15").Value = '  +0.14%  '
$ws.Range("E16").Value = '  -0.96%  '
$ws.Range("E17").Value = '  -1.05%  '
$ws.Range("D18").Value = '2.507.81'
$ws.Range("E18").Value = '  +1.50%  '
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("D20").Value = "'321.58"
$ws.Range("E20").Value = '  +0.50%  '
$ws.Range("D21").Value = "'4.17"
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").Value = "'0.997"
$ws.Range("E22").Value = '  -0.30%  '
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").Value = "'5.97"
$ws.Range("E23").Value = '  +4.46%  '
$ws.Range("D24").Value = "'64.21"
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("D25").Value = "'0.402"
$ws.Range("E25").Value = '  -1.92%  '
$ws.Range("D26").Value = "'0.994"
$ws.Range("E26").Value = '  -0.63%  '
$ws.Range("D28").Value = "'7.33"
$ws.Range("E28").Value = '  +0.10%  '
$ws.Range("D29").Value = '0.0₃0749'
$ws.Range("E29").Value = '  +0.14%  '
$ws.Range("D30").Value = "'167.96"
$ws.Range("E30").Value = '  +1.77%  '
$ws.Range("E31").Value = '  +1.18%  '
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").Value = "'1.18"
$ws.Range("E32").Value = '  +1.48%  '
$ws.Range("B33").Value = 'Aptos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D33").Value = "'6.24"
$ws.Range("E33").Value = '  +0.30%  '
$ws.Range("D34").Value = "'0.997"
$ws.Range("E34").Value = '  -0.12%  '
$ws.Range("D35").Value = "'0.993"
$ws.Range("E35").Value = '  -0.57%  '
$ws.Range("D36").Value = "'18.05"
$ws.Range("E36").Value = '  +0.20%  '
$ws.Range("E37").Value = '  -2.13%  '
$ws.Range("E38").Value = '  -0.98%  '
$ws.Range("D39").Value = "'36.83"
$ws.Range("E39").Value = '  +0.92%  '
$ws.Range("E40").Value = '  -0.85%  '
$ws.Range("D41").Value = "'0.769"
$ws.Range("E41").Value = '  -2.64%  '
$ws.Range("D42").Value = "'276.66"
$ws.Range("E42").Value = '  +1.81%  '
$ws.Range("E43").Value = '  +2.35%  '
$ws.Range("E44").Value = '  -0.21%  '
$ws.Range("D45").Value = "'0.596"
$ws.Range("E45").Value = '  +0.85%  '
$ws.Range("D46").Value = "'0.0920"
$ws.Range("E46").Value = '  +1.78%  '
$ws.Range("D47").Value = "'121.36"
$ws.Range("E47").Value = '  -4.20%  '
$ws.Range("E48").Value = '  +2.48%  '
$ws.Range("E49").Value = '  +0.18%  '
$ws.Range("E50").Value = '  +1.34%  '
$ws.Range("E51").Value = '  +0.01%  '

# Strip the quote-prefix style artifact so these cells keep the default
# (unstyled) cell format, matching the original workbook formatting -
# only the textual content changes.
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D12").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D47").ClearFormats()
